{"js": "// 1. The \"_GoBack\" bookmark that used to sit at the end of the paragraph\n//    ending in \"...for the data?\" is removed (it will be re-added at the\n//    end of the newly typed paragraph below, since that's where the\n//    author's cursor ended up after the edit).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Insert a brand new bullet (\"Currently, if you select a root directory...\")\n//    right after the paragraph that ends with \"...and the foreach loop.\"\n//    (the paragraph just above the \"To do:\" heading).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"foreach\") !== -1 && t.indexOf(\"loop.\") !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nconst newParagraph = anchorParagraph.insertParagraph(\"\", \"After\");\n// Give it the same numbered-list formatting as the rest of the \"Current bugs\" list.\nnewParagraph.attachToList(1, 0);\nawait context.sync();\n\n// Build the exact run/proofErr/bookmark structure for the new paragraph via OOXML,\n// since it is brand new content (nothing existing is being overwritten/collapsed).\nconst pkg =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Currently, if you select a root directory, such as drive C, D, E, </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>ect</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t>. The drive letter/name is not displayed for the output location.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nnewParagraph.insertOoxml(pkg, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. The \"_GoBack\" bookmark that used to sit at the end of the paragraph\n#    ending in \"...for the data?\" is removed (it will be re-added at the\n#    end of the newly typed paragraph below, since that's where the\n#    author's cursor ended up after the edit).\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# 2. Find the paragraph that ends with \"...and the foreach loop.\" (the\n#    paragraph just above the \"To do:\" heading) so we can insert the new\n#    bullet right after it.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"and the foreach loop.\") | Out-Null\n$findRange.Expand(4) | Out-Null   # wdParagraph\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $findRange.Start) {\n        $anchorIndex = $i\n        break\n    }\n}\n\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n\n# Build the exact run/proofErr/bookmark structure for the new paragraph via\n# WordOpenXML, since it is brand new content (nothing existing is being\n# overwritten/collapsed).\n$pkg = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Currently, if you select a root directory, such as drive C, D, E, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>ect</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>. The drive letter/name is not displayed for the output location.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$newParagraph.Range.InsertXML($pkg) | Out-Null\n"}
